# Generate Report for Handback
# Applies the "handback" report data to the Overview / zh-cn / de-de sheets:
#  - updates the status text from "Ready for handoff" to "Handed back: in sync with en-US"
#  - fills in the Latest Target File / Latest Handback File / Latest Handback DateTime
#    columns on the per-language sheets, with hyperlinks on the target-file column
#  - widens a couple of columns to fit the new, longer text

$wb = $excel.ActiveWorkbook

$ovr   = $wb.Worksheets.Item("Overview")
$zhcn  = $wb.Worksheets.Item("zh-cn")
$dede  = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the per-language status text, widen status columns
# ---------------------------------------------------------------------------
$ovr.Range("E2").Value = $newStatus
$ovr.Range("F2").Value = $newStatus
$ovr.Range("E3").Value = $newStatus
$ovr.Range("F3").Value = $newStatus

$ovr.Columns.Item(5).ColumnWidth = 29.14
$ovr.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("J2").Value = "12e9f1b5-392a-477a-b62d-c874087e6dd8.md"
$zhcn.Range("K2").Value = "12e9f1b5-392a-477a-b62d-c874087e6dd8.c3c1536ea5a3777d5ecdd81f83593c09a56e75d6.zh-cn.xlf"
$zhcn.Range("L2").Value = "2017-02-21 11:00:15"

$zhcn.Range("J3").Value = "cf5d4e63-38a4-456d-9f81-5cae228c4916.md"
$zhcn.Range("K3").Value = "cf5d4e63-38a4-456d-9f81-5cae228c4916.d0143481560cac4a9d527aece013bd34aa8d71cc.zh-cn.xlf"
$zhcn.Range("L3").Value = "2017-02-21 11:00:15"

$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(10).ColumnWidth = 39.17
$zhcn.Columns.Item(11).ColumnWidth = 39.17

# rebuild the hyperlinks in the desired order so relationship ids line up
# (A2, J2, A3, J3)
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/12e9f1b5-392a-477a-b62d-c874087e6dd8.md", $null, $null, "12e9f1b5-392a-477a-b62d-c874087e6dd8.md")
$zhcn.Hyperlinks.Add($zhcn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/12e9f1b5-392a-477a-b62d-c874087e6dd8.md", $null, $null, "12e9f1b5-392a-477a-b62d-c874087e6dd8.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/cf5d4e63-38a4-456d-9f81-5cae228c4916.md", $null, $null, "cf5d4e63-38a4-456d-9f81-5cae228c4916.md")
$zhcn.Hyperlinks.Add($zhcn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/cf5d4e63-38a4-456d-9f81-5cae228c4916.md", $null, $null, "cf5d4e63-38a4-456d-9f81-5cae228c4916.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("J2").Value = "12e9f1b5-392a-477a-b62d-c874087e6dd8.md"
$dede.Range("K2").Value = "12e9f1b5-392a-477a-b62d-c874087e6dd8.c3c1536ea5a3777d5ecdd81f83593c09a56e75d6.de-de.xlf"
$dede.Range("L2").Value = "2017-02-21 11:00:39"

$dede.Range("J3").Value = "cf5d4e63-38a4-456d-9f81-5cae228c4916.md"
$dede.Range("K3").Value = "cf5d4e63-38a4-456d-9f81-5cae228c4916.d0143481560cac4a9d527aece013bd34aa8d71cc.de-de.xlf"
$dede.Range("L3").Value = "2017-02-21 11:00:39"

$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(10).ColumnWidth = 39.17
$dede.Columns.Item(11).ColumnWidth = 39.17

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/12e9f1b5-392a-477a-b62d-c874087e6dd8.md", $null, $null, "12e9f1b5-392a-477a-b62d-c874087e6dd8.md")
$dede.Hyperlinks.Add($dede.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/12e9f1b5-392a-477a-b62d-c874087e6dd8.md", $null, $null, "12e9f1b5-392a-477a-b62d-c874087e6dd8.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/cf5d4e63-38a4-456d-9f81-5cae228c4916.md", $null, $null, "cf5d4e63-38a4-456d-9f81-5cae228c4916.md")
$dede.Hyperlinks.Add($dede.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8e0c226ac0de28608e52520b3b405d5b2b49027a/e2e/cf5d4e63-38a4-456d-9f81-5cae228c4916.md", $null, $null, "cf5d4e63-38a4-456d-9f81-5cae228c4916.md")

"Handback report generated"
